$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.614.85'

$ws.Range('E2').Value = '  +2.11%  '

$ws.Range('D3').Value = '1.887.89'

$ws.Range('E3').Value = '  +0.35%  '

$origStyle = $ws.Range('D4').Style
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = $origStyle

$ws.Range('E4').Value = '  +0.03%  '

$origStyle = $ws.Range('D5').Style
$ws.Range('D5').Value = "'245.15"
$ws.Range('D5').Style = $origStyle

$ws.Range('E5').Value = '  +0.94%  '

$ws.Range('E6').Value = '  +0.02%  '

$origStyle = $ws.Range('D7').Style
$ws.Range('D7').Value = "'0.4912"
$ws.Range('D7').Style = $origStyle

$ws.Range('E7').Value = '  -0.23%  '

$origStyle = $ws.Range('D8').Style
$ws.Range('D8').Value = "'0.2950"
$ws.Range('D8').Style = $origStyle

$ws.Range('E8').Value = '  +0.32%  '

$origStyle = $ws.Range('D9').Style
$ws.Range('D9').Value = "'0.06775"
$ws.Range('D9').Style = $origStyle

$ws.Range('E9').Value = '  +2.42%  '

$ws.Range('D10').Value = '1.887.67'

$ws.Range('E10').Value = '  +0.33%  '

$ws.Range('E11').Value = '  +3.61%  '

$origStyle = $ws.Range('D12').Style
$ws.Range('D12').Value = "'0.07234"
$ws.Range('D12').Style = $origStyle

$ws.Range('E12').Value = '  +0.73%  '

$origStyle = $ws.Range('D13').Style
$ws.Range('D13').Value = "'91.17"
$ws.Range('D13').Style = $origStyle

$ws.Range('E13').Value = '  +5.59%  '

$origStyle = $ws.Range('D14').Style
$ws.Range('D14').Value = "'0.6782"
$ws.Range('D14').Style = $origStyle

$ws.Range('E14').Value = '  +1.78%  '

$origStyle = $ws.Range('D15').Style
$ws.Range('D15').Value = "'5.042"
$ws.Range('D15').Style = $origStyle

$ws.Range('E15').Value = '  +3.53%  '

$ws.Range('D16').Value = '30.599.57'

$ws.Range('E16').Value = '  +2.12%  '

$origStyle = $ws.Range('D17').Style
$ws.Range('D17').Value = "'0.000007979"
$ws.Range('D17').Style = $origStyle

$ws.Range('E17').Value = '  +2.22%  '

$origStyle = $ws.Range('D18').Style
$ws.Range('D18').Value = "'1.000"
$ws.Range('D18').Style = $origStyle

$ws.Range('E18').Value = '  +0.03%  '

$origStyle = $ws.Range('D19').Style
$ws.Range('D19').Value = "'13.15"
$ws.Range('D19').Style = $origStyle

$ws.Range('E19').Value = '  +2.93%  '

$ws.Range('D20').Value = '2.131.46'

$ws.Range('E20').Value = '  +0.38%  '

$ws.Range('E21').Value = '  +0.16%  '

$origStyle = $ws.Range('D22').Style
$ws.Range('D22').Value = "'4.822"
$ws.Range('D22').Style = $origStyle

$ws.Range('E22').Value = '  +0.90%  '

$origStyle = $ws.Range('D23').Style
$ws.Range('D23').Value = "'193.35"
$ws.Range('D23').Style = $origStyle

$ws.Range('E23').Value = '  +37.26%  '

$origStyle = $ws.Range('D24').Style
$ws.Range('D24').Value = "'6.072"
$ws.Range('D24').Style = $origStyle

$ws.Range('E24').Value = '  +3.93%  '

$origStyle = $ws.Range('D25').Style
$ws.Range('D25').Value = "'9.328"
$ws.Range('D25').Style = $origStyle

$ws.Range('E25').Value = '  +2.76%  '

$origStyle = $ws.Range('D26').Style
$ws.Range('D26').Value = "'155.49"
$ws.Range('D26').Style = $origStyle

$ws.Range('E26').Value = '  +3.40%  '

$origStyle = $ws.Range('D27').Style
$ws.Range('D27').Value = "'19.20"
$ws.Range('D27').Style = $origStyle

$ws.Range('E27').Value = '  +13.36%  '

$origStyle = $ws.Range('D28').Style
$ws.Range('D28').Value = "'1.904"
$ws.Range('D28').Style = $origStyle

$ws.Range('E28').Value = '  +0.33%  '

$origStyle = $ws.Range('D29').Style
$ws.Range('D29').Value = "'1.399"
$ws.Range('D29').Style = $origStyle

$ws.Range('E29').Value = '  +0.53%  '

$origStyle = $ws.Range('D30').Style
$ws.Range('D30').Value = "'4.327"
$ws.Range('D30').Style = $origStyle

$ws.Range('E30').Value = '  +3.21%  '

$origStyle = $ws.Range('D31').Style
$ws.Range('D31').Value = "'0.09056"
$ws.Range('D31').Style = $origStyle

$ws.Range('E31').Value = '  +3.44%  '

$origStyle = $ws.Range('D32').Style
$ws.Range('D32').Value = "'4.021"
$ws.Range('D32').Style = $origStyle

$ws.Range('E32').Value = '  +1.02%  '

$origStyle = $ws.Range('D33').Style
$ws.Range('D33').Value = "'0.05202"
$ws.Range('D33').Style = $origStyle

$ws.Range('E33').Value = '  +3.90%  '

$origStyle = $ws.Range('D34').Style
$ws.Range('D34').Value = "'0.7550"
$ws.Range('D34').Style = $origStyle

$ws.Range('E34').Value = '  +5.12%  '

$origStyle = $ws.Range('D35').Style
$ws.Range('D35').Value = "'1.111"
$ws.Range('D35').Style = $origStyle

$ws.Range('E35').Value = '  +0.06%  '

$origStyle = $ws.Range('D36').Style
$ws.Range('D36').Value = "'2.759"
$ws.Range('D36').Style = $origStyle

$ws.Range('E36').Value = '  +3.38%  '

$origStyle = $ws.Range('D37').Style
$ws.Range('D37').Value = "'0.01836"
$ws.Range('D37').Style = $origStyle

$ws.Range('E37').Value = '  +2.64%  '

$ws.Range('E38').Value = '  -1.01%  '

$origStyle = $ws.Range('D39').Style
$ws.Range('D39').Value = "'2.145"
$ws.Range('D39').Style = $origStyle

$ws.Range('E39').Value = '  -0.62%  '

$origStyle = $ws.Range('D40').Style
$ws.Range('D40').Value = "'0.9371"
$ws.Range('D40').Style = $origStyle

$ws.Range('E40').Value = '  -0.36%  '

$origStyle = $ws.Range('D41').Style
$ws.Range('D41').Value = "'0.4418"
$ws.Range('D41').Style = $origStyle

$ws.Range('E41').Value = '  +4.53%  '

$origStyle = $ws.Range('D42').Style
$ws.Range('D42').Value = "'105.10"
$ws.Range('D42').Style = $origStyle

$ws.Range('E42').Value = '  +1.47%  '

$ws.Range('E43').Value = '  +0.07%  '

$origStyle = $ws.Range('D44').Style
$ws.Range('D44').Value = "'5.739"
$ws.Range('D44').Style = $origStyle

$ws.Range('E44').Value = '  +0.07%  '

$origStyle = $ws.Range('D45').Style
$ws.Range('D45').Value = "'7.598"
$ws.Range('D45').Style = $origStyle

$ws.Range('E45').Value = '  +3.71%  '

$origStyle = $ws.Range('D46').Style
$ws.Range('D46').Value = "'0.1346"
$ws.Range('D46').Style = $origStyle

$ws.Range('E46').Value = '  +6.17%  '

$origStyle = $ws.Range('D47').Style
$ws.Range('D47').Value = "'0.05856"
$ws.Range('D47').Style = $origStyle

$ws.Range('E47').Value = '  +2.78%  '

$ws.Range('B48').Value = 'EnergySwap'

$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$origStyle = $ws.Range('D48').Style
$ws.Range('D48').Value = "'8.711"
$ws.Range('D48').Style = $origStyle

$ws.Range('E48').Value = '  +5.43%  '

$ws.Range('B49').Value = 'NEARProtocol'

$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$origStyle = $ws.Range('D49').Style
$ws.Range('D49').Value = "'1.432"
$ws.Range('D49').Style = $origStyle

$ws.Range('E49').Value = '  +7.03%  '

$origStyle = $ws.Range('D50').Style
$ws.Range('D50').Value = "'0.3921"
$ws.Range('D50').Style = $origStyle

$ws.Range('E50').Value = '  +4.41%  '

$ws.Range('E51').Value = '  +2.53%  '
